$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.209.51"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.028.00"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0788"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.34%  "
$ws.Range("D12").Value = "2.317.14"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.743"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "2.033.64"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "37.153.33"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "0.0₃0817"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("E25").Value = "  -4.75%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.59%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("E35").Value = "  -4.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  +6.49%  "
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").Value = "1.473.19"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0915"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("E46").Value = "  -5.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "2.210.59"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -8.11%  "
